$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.480.22'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.70%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.974.72'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.91%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.31%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '326.74'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.37%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.003'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.22%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4658'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.19%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3916'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.23%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '46.17'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.29%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07951'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.79%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9917'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.20%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.82'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.71%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.973.52'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.29%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.179'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.33%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.850'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.71%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.07085'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.46%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '87.64'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.88%  '

# Row 18
$ws.Range('E18').Value = '  +0.30%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000009943'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.38%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.29'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.10%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.004'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.38%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '29.479.35'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.72%  '

# Row 23
$ws.Range('B23').Value = 'BitDAO'
$ws.Range('C23').Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.5031'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.33%  '

# Row 24
$ws.Range('B24').Value = 'Uniswap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.534'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.05%  '

# Row 25
$ws.Range('B25').Value = 'Cosmos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.16'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.39%  '

# Row 26
$ws.Range('B26').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C26').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.212.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.58%  '

# Row 27
$ws.Range('B27').Value = 'Toncoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.107'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.42%  '

# Row 28
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '158.50'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.78%  '

# Row 29
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.54'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.32%  '

# Row 30
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.798'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.94%  '

# Row 31
$ws.Range('B31').Value = 'BitcoinCash'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '119.62'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.71%  '

# Row 32
$ws.Range('B32').Value = 'LidoDAOToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.907'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.14%  '

# Row 33
$ws.Range('B33').Value = 'Stellar'
$ws.Range('C33').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09414'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.46%  '

# Row 34
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.8921'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.98%  '

# Row 35
$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.232'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.05%  '

# Row 36
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.322'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.37%  '

# Row 37
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.172'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.41%  '

# Row 38
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05818'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.47%  '

# Row 39
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.171'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.53%  '

# Row 40
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.02105'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.59%  '

# Row 41
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.765'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.43%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.000003224'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +47.75%  '

# Row 43
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.5715'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.07%  '

# Row 44
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1799'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.74%  '

# Row 45
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '9.662'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.10%  '

# Row 46
$ws.Range('B46').Value = 'MXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.755'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.96%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '11.82'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.32%  '

# Row 48
$ws.Range('B48').Value = 'Decentraland'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5356'
$ws.Range('D48').Style = 'Normal'

# Row 49
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.195'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.02%  '

# Row 50
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06924'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.58%  '

# Row 51
$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '114.14'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.45%  '
